$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Phuc123123"
$ws.Range("B3").Value = "123456aA@"
